$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 32466.666
$ws.Range("J3").Value = 32466.666
$ws.Range("L3").Value = 32466.666
$ws.Range("N3").Value = -32694.666
$ws.Range("H43").Value = 834.94116
$ws.Range("I43").Value = 667
$ws.Range("J43").Value = 870.9286
$ws.Range("K43").Value = 667
$ws.Range("L43").Value = 870.9286
$ws.Range("M43").Value = -598
$ws.Range("N43").Value = -1008.9286
$ws.Range("H102").Value = 32466.666
$ws.Range("J102").Value = 32466.666
$ws.Range("L102").Value = 32466.666
$ws.Range("N102").Value = -38956.666
$ws.Range("H111").Value = 1916.6666
$ws.Range("I111").Value = 1916.6666
$ws.Range("K111").Value = 5749.9998
$ws.Range("M111").Value = -2682.9998
$ws.Range("H127").Value = 824.875
$ws.Range("I127").Value = 232.25
$ws.Range("J127").Value = 1417.5
$ws.Range("K127").Value = 696.75
$ws.Range("L127").Value = 4252.5
$ws.Range("M127").Value = 4263.25
$ws.Range("N127").Value = -14172.5
$ws.Range("H132").Value = 6536.5
$ws.Range("I132").Value = 5598.5
$ws.Range("J132").Value = 8814.5
$ws.Range("K132").Value = 16795.5
$ws.Range("L132").Value = 26443.5
$ws.Range("M132").Value = -14265.5
$ws.Range("N132").Value = -31503.5

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1641819.1
$ws.Range("I32").Value = 11448.218
$ws.Range("J32").Value = 7562639.5
$ws.Range("K32").Value = 11448.218
$ws.Range("L32").Value = 7562639.5
$ws.Range("M32").Value = -11161.218
$ws.Range("N32").Value = -7563213.5
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H122").Value = 2090.3635
$ws.Range("I122").Value = 3600
$ws.Range("K122").Value = 10800
$ws.Range("M122").Value = -8350
$ws.Range("H132").Value = 2278292
$ws.Range("I132").Value = 5589.6924
$ws.Range("J132").Value = 4315887.5
$ws.Range("K132").Value = 16769.0772
$ws.Range("L132").Value = 12947662.5
$ws.Range("M132").Value = -14239.0772
$ws.Range("N132").Value = -12952722.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 850.6
$ws.Range("I22").Value = 320.8
$ws.Range("J22").Value = 1380.4
$ws.Range("K22").Value = 320.8
$ws.Range("L22").Value = 1380.4
$ws.Range("M22").Value = 29.19999999999999
$ws.Range("N22").Value = -2080.4
$ws.Range("H43").Value = 21282.125
$ws.Range("J43").Value = 21282.125
$ws.Range("L43").Value = 21282.125
$ws.Range("N43").Value = -21650.125
$ws.Range("H101").Value = 21282.125
$ws.Range("J101").Value = 21282.125
$ws.Range("L101").Value = 21282.125
$ws.Range("N101").Value = -27772.125
$ws.Range("H132").Value = 2136.1052
$ws.Range("I132").Value = 1552.45
$ws.Range("J132").Value = 2784.611
$ws.Range("K132").Value = 4657.35
$ws.Range("L132").Value = 8353.832999999999
$ws.Range("M132").Value = -2127.35
$ws.Range("N132").Value = -13413.833

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 703.63635
$ws.Range("I113").Value = 640.7273
$ws.Range("J113").Value = 766.5454999999999
$ws.Range("K113").Value = 1922.1819
$ws.Range("L113").Value = 2299.6365
$ws.Range("M113").Value = 247.8181
$ws.Range("N113").Value = -6639.6365
$ws.Range("H136").Value = 3230
$ws.Range("I136").Value = 3230
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9690
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4590
$ws.Range("N136").ClearContents()

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -471
$ws.Range("H25").Value = 1500
$ws.Range("J25").Value = 1500
$ws.Range("L25").Value = 1500
$ws.Range("N25").Value = -2558
$ws.Range("H43").Value = 12500
$ws.Range("I43").Value = 6000
$ws.Range("J43").Value = 19000
$ws.Range("K43").Value = 6000
$ws.Range("L43").Value = 19000
$ws.Range("M43").Value = -5849
$ws.Range("N43").Value = -19302
$ws.Range("H70").Value = 7935.5557
$ws.Range("I70").Value = 9387.5
$ws.Range("K70").Value = 9387.5
$ws.Range("M70").Value = -9117.5
$ws.Range("H73").Value = 7935.5557
$ws.Range("I73").Value = 9387.5
$ws.Range("K73").Value = 9387.5
$ws.Range("M73").Value = -8451.5
$ws.Range("H74").Value = 13000
$ws.Range("J74").Value = 13000
$ws.Range("L74").Value = 13000
$ws.Range("N74").Value = -14872
$ws.Range("H75").Value = 2000
$ws.Range("I75").Value = 2000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 2000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -1126
$ws.Range("N75").ClearContents()
$ws.Range("H77").Value = 13000
$ws.Range("J77").Value = 13000
$ws.Range("L77").Value = 39000
$ws.Range("N77").Value = -48360
$ws.Range("H78").Value = 2000
$ws.Range("I78").Value = 2000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 6000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -1632
$ws.Range("N78").ClearContents()
$ws.Range("H80").Value = 3000.75
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3002
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3002
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -4998
$ws.Range("H83").Value = 3000.75
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3002
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 15010
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -24994
$ws.Range("H122").Value = 3799.4285
$ws.Range("I122").Value = 3919.2
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 11757.6
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -9307.599999999999
$ws.Range("N122").Value = -15400
$ws.Range("H123").Value = 19732.5
$ws.Range("J123").Value = 19732.5
$ws.Range("L123").Value = 19732.5
$ws.Range("N123").Value = -24632.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 358.125
$ws.Range("I9").Value = 227.5
$ws.Range("J9").Value = 750
$ws.Range("K9").Value = 227.5
$ws.Range("L9").Value = 750
$ws.Range("M9").Value = -3.5
$ws.Range("N9").Value = -1198
$ws.Range("H103").Value = 51650.5
$ws.Range("J103").Value = 51650.5
$ws.Range("L103").Value = 51650.5
$ws.Range("N103").Value = -53994.5
$ws.Range("H132").Value = 15828
$ws.Range("I132").Value = 18269.357
$ws.Range("J132").Value = 7219
$ws.Range("K132").Value = 54808.071
$ws.Range("L132").Value = 21657
$ws.Range("M132").Value = -52278.071
$ws.Range("N132").Value = -26717

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2798.5557
$ws.Range("I122").Value = 2708.3333
$ws.Range("J122").Value = 2888.7778
$ws.Range("K122").Value = 8124.999899999999
$ws.Range("L122").Value = 8666.3334
$ws.Range("M122").Value = -5674.999899999999
$ws.Range("N122").Value = -13566.3334
$ws.Range("H126").Value = 2689.5789
$ws.Range("I126").Value = 3110.2
$ws.Range("J126").Value = 2222.2222
$ws.Range("K126").Value = 9330.599999999999
$ws.Range("L126").Value = 6666.6666
$ws.Range("M126").Value = -6860.599999999999
$ws.Range("N126").Value = -11606.6666
$ws.Range("H132").Value = 1982.3864
$ws.Range("I132").Value = 1632.9131
$ws.Range("J132").Value = 2365.1428
$ws.Range("K132").Value = 4898.7393
$ws.Range("L132").Value = 7095.428400000001
$ws.Range("M132").Value = -2368.7393
$ws.Range("N132").Value = -12155.4284
